$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# The "Status" text is shared by the Overview summary (zh-cn / de-de columns)
# and each language sheet's own Status column - update them all so the
# report reflects the handback being in sync with en-US.
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Refresh the Latest Handback DateTime stamps for both languages
$wsZhCn.Range("K2").Value = "2016-10-14 07:57:14"
$wsDeDe.Range("K2").Value = "2016-10-14 07:57:31"

# The handback is now current, so the stale "not the latest" error goes away
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Resize the affected columns to match the regenerated report widths
$wsOverview.Range("E:F").ColumnWidth = 29.1666666666667
$wsZhCn.Range("C:C").ColumnWidth = 29.1666666666667
$wsZhCn.Range("P:P").ColumnWidth = 12.8333333333333
$wsDeDe.Range("C:C").ColumnWidth = 29.1666666666667
$wsDeDe.Range("P:P").ColumnWidth = 12.8333333333333
